# Applies the "Add files via upload" edit to the workbook:
#  - Sheet "apartado8" (sheet1): adds LaTeX-style helper formulas in columns
#    E/F for rows 2-9 (a stray fill-down leaves one extra shared formula in
#    F10), and relocates the small "fitting in an exponential" summary block
#    that used to sit at F4:H8 down to B13:D17.
#  - Sheet "apartado4" (sheet2): only the active selection changes.
#  - Sheet "apartado5" (sheet3): adds the same style of E/F helper formulas
#    for rows 2-7.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: apartado8
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("apartado8")

# Relocate the "fitting in an exponential" block from F4:H8 to B13:D17
# before adding the new formulas, so nothing overlaps the E:F columns.
$ws1.Range("B13").Value2 = $ws1.Range("F4").Value2

$ws1.Range("B14").Value2 = $ws1.Range("F5").Value2
$ws1.Range("C14").Value2 = $ws1.Range("G5").Value2
$ws1.Range("D14").Value2 = $ws1.Range("H5").Value2

$ws1.Range("B15").Value2 = $ws1.Range("F6").Value2
$ws1.Range("C15").Value2 = $ws1.Range("G6").Value2
$ws1.Range("D15").Value2 = $ws1.Range("H6").Value2

$ws1.Range("B17").Value2 = $ws1.Range("F8").Value2
$ws1.Range("C17").Value2 = $ws1.Range("G8").Value2

$ws1.Range("F4:H8").ClearContents() | Out-Null

# New helper columns: E -> "$N$", F -> "$K\pm errK$"
$ws1.Range("E2").Formula = '="$" & A2 & "$"'
$ws1.Range("F2").Formula = '="$" & B2 & "\pm" & C2 & "$"'

# Fill down E3:E9 / F3:F10 with relative references (F was dragged one row
# too far by the original author, producing the stray F10 "$\pm$" cell).
$ws1.Range("E3:E9").Formula = '="$" & A3 & "$"'
$ws1.Range("F3:F10").Formula = '="$" & B3 & "\pm" & C3 & "$"'

$ws1.Range("E2:F9").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: apartado4 (only the selection moves)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("apartado4")
$ws2.Activate() | Out-Null
$ws2.Range("C19").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: apartado5
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("apartado5")

$ws3.Range("E2").Formula = '="$" & A2 & "$"'
$ws3.Range("F2").Formula = '="$" & B2 & "\pm" & C2 & "$"'

$ws3.Range("E3:E7").Formula = '="$" & A3 & "$"'
$ws3.Range("F3:F7").Formula = '="$" & B3 & "\pm" & C3 & "$"'

$ws3.Activate() | Out-Null
$ws3.Range("E2:F7").Select() | Out-Null

# ---------------------------------------------------------------------
# Leave sheet1 active/selected, as in the target workbook.
# ---------------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("E2:F9").Select() | Out-Null
